# Reproduces the "week 5 assignment submitted" edit:
#   - renumbers/extends the documents hyperlink relationships by adding 4 new
#     reference hyperlinks (the existing 2 are left in place; Word renumbers
#     relationship ids for us on save)
#   - inserts a whole new discussion-reply block (quote + response + References
#     list + a second quote + reply) ahead of the documents final "Emanuel"
#     sign-off paragraph, which becomes the sign-off for the new content

$d = $word.ActiveDocument

# Insert all of the new paragraphs (as plain text, `r = new Word paragraph)
# immediately before the very last paragraph ("Emanuel"), which is left alone.
$anchorPara = $d.Paragraphs.Last
$ins = $anchorPara.Range.Duplicate
$ins.Collapse(1)
$newText = "`rEmanuel`r`r`r“Regarding malware, if we compare windows and linux OS which do you believe is the most affect by malware?”`r`rSince approximately 90% of all desktop operating systems in use are Windows while Linux has a paltry 1.6% (Net Applications, 2014), it doesn’t come as a surprise that it is the most attacked OS (Vaughan-Nichols, 2014). Since hackers “want to have a higher success rate”, they will of course target Windows (Popa, 2014). However, it has been recently reported that approximately ”25,000 web servers infected with Linux malware have been used in the past two years to hit website visitors” (Tung, 2014). It just goes to show that despite Linux’s reputation as being secure, it really isn’t.`r`rReferences:`r`rNet Applications (2014) Desktop Operating System Market Share [Online]. Available from: http://www.netmarketshare.com/operating-system-market-share.aspx?qprid=10&qpcustomd=0 (Accessed: 18 June 2014)`r`rPopa,  B. (2014) ‘Security Expert Explains Why Windows Is the Most Attacked Operating System’, Softpedia [Online]. Available from: http://news.softpedia.com/news/Security-Expert-Explains-Why-Windows-Is-the-Most-Attacked-Operating-System-445834.shtml (Accessed: 18 June 2014)`r`rTung, L. (2014) ‘Botnet of thousands of Linux servers pumps Windows desktop malware onto web’, ZDNet [Online]. Available from: http://www.zdnet.com/botnet-of-thousands-of-linux-servers-pumps-windows-desktop-malware-onto-web-7000027472/ (Accessed: 18 June 2014)`r`rVaughan-Nichols, S. (2014) ‘Security 2014: The holes are in the apps, not the operating systems’, ZDNet [Online]. Available from: http://www.zdnet.com/security-2014-the-holes-are-in-the-apps-not-the-operating-systems-7000026893/ (Accessed: 18 June 2014)`r`r`r“The more secure an application the less user friendly it is. So we have to think about what is important. Is making a very beautiful and user friendly application that provides a rich user experience more important that securing your application, user dat? In my eyes no matter what, making sure your application is robust and secure is more important than the user experience.”`r`rHi Adrian,`r`rThat’s an interesting point. I’ve always thought that user-friendliness and power/flexibility were inversely proportional, but never thought of it also applying to user-friendliness and security. Regarding which is more important, I think it would depend on the customer. Individual consumers will probably have less stringent security requirements, while businesses and larger organizations may have stricter ones. For example, everyone is required to change their passwords every 90 days where I work. For my own personal accounts (email, banking, etc.), no such requirements are in place, though it probably is a good idea. I have so many of these, that requirements to change them periodically would be obtrusive. However, knowing that events such as irregular activity will result in someone calling me offers me some peace of mind.`r`rRegards,`r`r"
$ins.InsertAfter($newText)

# Italicise the "Softpedia" source-title run.
$rng = $d.Content
$null = $rng.Find.Execute("Softpedia ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Italic = $true

# Turn the 4 new reference URLs into real hyperlinks.
$rng = $d.Content
$null = $rng.Find.Execute("http://www.netmarketshare.com/operating-system-market-share.aspx?qprid=10&qpcustomd=0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$null = $d.Hyperlinks.Add($rng, "http://www.netmarketshare.com/operating-system-market-share.aspx?qprid=10&qpcustomd=0", "", "", "http://www.netmarketshare.com/operating-system-market-share.aspx?qprid=10&qpcustomd=0")
$rng = $d.Content
$null = $rng.Find.Execute("http://news.softpedia.com/news/Security-Expert-Explains-Why-Windows-Is-the-Most-Attacked-Operating-System-445834.shtml", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$null = $d.Hyperlinks.Add($rng, "http://news.softpedia.com/news/Security-Expert-Explains-Why-Windows-Is-the-Most-Attacked-Operating-System-445834.shtml", "", "", "http://news.softpedia.com/news/Security-Expert-Explains-Why-Windows-Is-the-Most-Attacked-Operating-System-445834.shtml")
$rng = $d.Content
$null = $rng.Find.Execute("http://www.zdnet.com/botnet-of-thousands-of-linux-servers-pumps-windows-desktop-malware-onto-web-7000027472/", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$null = $d.Hyperlinks.Add($rng, "http://www.zdnet.com/botnet-of-thousands-of-linux-servers-pumps-windows-desktop-malware-onto-web-7000027472/", "", "", "http://www.zdnet.com/botnet-of-thousands-of-linux-servers-pumps-windows-desktop-malware-onto-web-7000027472/")
$rng = $d.Content
$null = $rng.Find.Execute("http://www.zdnet.com/security-2014-the-holes-are-in-the-apps-not-the-operating-systems-7000026893/", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$null = $d.Hyperlinks.Add($rng, "http://www.zdnet.com/security-2014-the-holes-are-in-the-apps-not-the-operating-systems-7000026893/", "", "", "http://www.zdnet.com/security-2014-the-holes-are-in-the-apps-not-the-operating-systems-7000026893/")

Write-Output ("paragraphs=" + $d.Paragraphs.Count + " hyperlinks=" + $d.Hyperlinks.Count)
